# lipidcane_spearman_2.xlsx -- further work on lipidcane2g
# - Rename several metric / parameter headers (ton-basis -> yearly-basis units)
# - Add a new "Heat exchanger network error [%]" metric column (J)
# - Re-order / rename several parameter rows (Lipid content / Lipid retention /
#   Additional lipid extraction efficiency / Capacity)
# - Drop the old "Fermentation" / "Solids loading [%]" row entirely
# - Refresh every Spearman correlation value in the table

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Drop row 14 (Fermentation / Solids loading [%]) entirely, shifting rows up
# ---------------------------------------------------------------------------
$ws.Rows(14).Delete()

# ---------------------------------------------------------------------------
# 2) Metric header row (row 2) -- relabel to yearly-production units
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = "Biodiesel production [MMGal/yr]"
$ws.Range("E2").Value = "Ethanol production [MMGal/yr]"
$ws.Range("F2").Value = "Electricity production [MMWhr/yr]"
$ws.Range("G2").Value = "Natural gas consumption [MMcf/yr]"
$ws.Range("H2").Value = "Productivity [MMGGE/yr]"
# I2 "TCI [10^6*USD]" is unchanged

# ---------------------------------------------------------------------------
# 3) Add the new metric column J: "Heat exchanger network error [%]"
# ---------------------------------------------------------------------------
$ws.Range("J2").Value = "Heat exchanger network error [%]"
$ws.Range("J2").Font.Bold = $true
$ws.Range("J2").HorizontalAlignment = -4108
$ws.Range("J2").VerticalAlignment = -4160
$ws.Range("J2").Borders.LineStyle = 1

# Extend the title merge C1:I1 -> C1:J1
$ws.Range("C1:J1").MergeCells = $true
$ws.Range("C1:J1").Borders.LineStyle = 1

# J1 needs the same header style as the rest of row 1
$ws.Range("B1:J1").Font.Bold = $true
$ws.Range("B1:J1").HorizontalAlignment = -4108
$ws.Range("B1:J1").VerticalAlignment = -4160
$ws.Range("B1:J1").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 4) Parameter rows -- relabel the lipidcane stream's parameters and reorder
#    (old: Capacity, Lipid content, Efficiency, Lipid retention)
#    (new: Lipid content, Lipid retention, Additional lipid extraction
#          efficiency, Capacity)
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "Lipid content [dry wt. %]"
$ws.Range("B5").Value = "Lipid retention [%]"
$ws.Range("B6").Value = "Additional lipid extraction efficiency [%]"
$ws.Range("B7").Value = "Capacity [ton/hr]"

# ---------------------------------------------------------------------------
# 5) Refresh every Spearman correlation value (columns C, D, E, G, H, I, J)
#    for rows 4-13
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = 0.1529719852548794
$ws.Range("D4").Value = 0.9692952417478096
$ws.Range("E4").Value = -0.6891081906363276
$ws.Range("G4").Value = -0.8777052254122089
$ws.Range("H4").Value = 0.5872650572826023
$ws.Range("I4").Value = 0.3733875746475029
$ws.Range("J4").Value = 0.8363136114667618

$ws.Range("C5").Value = -0.01558384075135363
$ws.Range("D5").Value = 0.0001386908215476328
$ws.Range("E5").Value = 0.004659716730388668
$ws.Range("G5").Value = 0.004156849414273977
$ws.Range("H5").Value = 0.001202094576083783
$ws.Range("I5").Value = 0.01586069861842794
$ws.Range("J5").Value = 0.004601810299496158

$ws.Range("C6").Value = 0.03625086644203465
$ws.Range("D6").Value = 0.04710816447632658
$ws.Range("E6").Value = -0.002926963029078521
$ws.Range("G6").Value = 0.05862890634515625
$ws.Range("H6").Value = 0.01933925798957032
$ws.Range("I6").Value = 0.007594343631773744
$ws.Range("J6").Value = -0.01087775847520339

$ws.Range("C7").Value = 0.0804326726893069
$ws.Range("D7").Value = 0.1706032100241284
$ws.Range("E7").Value = 0.5662937319317493
$ws.Range("G7").Value = -0.09132257750890309
$ws.Range("H7").Value = 0.7044170071206801
$ws.Range("I7").Value = 0.9179981095199241
$ws.Range("J7").Value = 0.129208128128456

$ws.Range("C8").Value = 0.8479121154524845
$ws.Range("D8").Value = 0.0005414165976566638
$ws.Range("E8").Value = -0.02024814954592598
$ws.Range("G8").Value = -0.01129575299583012
$ws.Range("H8").Value = -0.01771302445252098
$ws.Range("I8").Value = -0.02237395174295807
$ws.Range("J8").Value = 0.006662663695463571

$ws.Range("C9").Value = 0.3355343397893736
$ws.Range("D9").Value = -0.008301011948040476
$ws.Range("E9").Value = 0.01556498606259944
$ws.Range("G9").Value = 0.007315422916616915
$ws.Range("H9").Value = -0.0002317522652700906
$ws.Range("I9").Value = 0.001721262884850515
$ws.Range("J9").Value = -0.0102886266111998

$ws.Range("C10").Value = -0.008654710522188419
$ws.Range("D10").Value = 0.006022146672885866
$ws.Range("E10").Value = -0.0106797724271909
$ws.Range("G10").Value = -0.0008984959079398362
$ws.Range("H10").Value = -0.006352999742119988
$ws.Range("I10").Value = -0.01553317540532701
$ws.Range("J10").Value = -0.009339791804540476

$ws.Range("C11").Value = -0.03648181797127271
$ws.Range("D11").Value = 0.0002448177697927107
$ws.Range("E11").Value = -0.002838362705534508
$ws.Range("G11").Value = -0.003835427577417102
$ws.Range("H11").Value = -0.0007139766045590642
$ws.Range("I11").Value = 0.00229044076361763
$ws.Range("J11").Value = 0.01858180405309702

$ws.Range("C12").Value = 0.1265808097672324
$ws.Range("D12").Value = 0.1646865822354633
$ws.Range("E12").Value = 0.4111203164448126
$ws.Range("G12").Value = 0.3875048927481957
$ws.Range("H12").Value = 0.3607830423673216
$ws.Range("I12").Value = 0.01017872699914908
$ws.Range("J12").Value = 0.02245004697870471

$ws.Range("C13").Value = -0.2608760032990401
$ws.Range("D13").Value = 0.01005232849809314
$ws.Range("E13").Value = -0.01249742142789685
$ws.Range("G13").Value = -0.01687056912282276
$ws.Range("H13").Value = -0.002855718738228749
$ws.Range("I13").Value = 0.005119967628798705
$ws.Range("J13").Value = 0.009960516405831281
